# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the outline properties used on the other sheets (summaryBelow/summaryRight)
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.Outline.SummaryRow = 1

# Copy the header formatting (bold + border) from an existing header cell
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Copy the date-column formatting from an existing date cell
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)  # xlPasteFormats

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$data = @(
    @(45550.99999999999, 308, 207.470108502219,  406.5618481621268),
    @(45578.99999999999, 186, 91.04464182403093, 286.9842738367208),
    @(45585.99999999999, 155, 56.81096716924829, 249.2743703609488),
    @(45599.99999999999, 94,  -4.414251926061912, 198.178165880847),
    @(45606.99999999999, 63,  -43.03810540953737, 159.6883095510039),
    @(45613.99999999999, 33,  -65.48310016526452, 129.498885292139),
    @(45620.99999999999, 2,   -103.6221863063046, 104.551879402687),
    @(45627.99999999999, 0,   -126.1452682058185, 70.25410627111035),
    @(45634.99999999999, 0,   -154.1591738348389, 40.45448145463535),
    @(45641.99999999999, 0,   -185.9573520000343, 10.43251764506441),
    @(45648.99999999999, 0,   -221.756744366632,  -20.41413934106603),
    @(45655.99999999999, 0,   -251.6459100340288, -48.99598281908152)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsForecast.Range("A1").Select()
